# Generate Report for Handoff
# Adds two new tracked files (d4eb6e29-... and f1d821ef-...) to the
# localization-status report: one new row on "Overview", and one new
# row on each of the "zh-cn" / "de-de" language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------

# Row 4 - d4eb6e29-4a72-4046-9332-9f9e7f4f6245
$overview.Range("A4").Value2 = "d4eb6e29-4a72-4046-9332-9f9e7f4f6245.md"
$overview.Hyperlinks.Add($overview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/d4eb6e29-4a72-4046-9332-9f9e7f4f6245.md", "", "", "d4eb6e29-4a72-4046-9332-9f9e7f4f6245.md") | Out-Null
$overview.Range("B4").Value2 = "Ready for handoff"
$overview.Range("C4").Value2 = "Ready for handoff"
$overview.Range("D4").Value2 = "2016-03-21 14:37:58"
$overview.Range("D4").NumberFormat = $dateFmt

# Row 5 - f1d821ef-a8ce-4150-9380-8d83689f79af
$overview.Range("A5").Value2 = "f1d821ef-a8ce-4150-9380-8d83689f79af.md"
$overview.Hyperlinks.Add($overview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/f1d821ef-a8ce-4150-9380-8d83689f79af.md", "", "", "f1d821ef-a8ce-4150-9380-8d83689f79af.md") | Out-Null
$overview.Range("B5").Value2 = "Ready for handoff"
$overview.Range("C5").Value2 = "Ready for handoff"
$overview.Range("D5").Value2 = "2016-03-21 14:37:58"
$overview.Range("D5").NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Helper data for the per-language detail sheets (zh-cn / de-de)
# Columns: A Source File Name | B File Extension | C Status |
#          D Latest Handoff File | E Latest Handoff Datetime |
#          F Latest Target File | G Latest Handback File |
#          H Latest Handback DateTime | I Reference Tokens |
#          J Handoff Reason | K Dependency From | L Error Detail
# ---------------------------------------------------------------------

function Add-DetailRows($ws, $langSuffix, $handoffDatetime) {
    # Row 4 - d4eb6e29-4a72-4046-9332-9f9e7f4f6245
    $ws.Range("A4").Value2 = "d4eb6e29-4a72-4046-9332-9f9e7f4f6245.md"
    $ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.$langSuffix/blob/master/e2e/d4eb6e29-4a72-4046-9332-9f9e7f4f6245.md", "", "", "d4eb6e29-4a72-4046-9332-9f9e7f4f6245.md") | Out-Null
    $ws.Range("B4").Value2 = ".md"
    $ws.Range("C4").Value2 = "Ready for handoff"
    $targetFile4 = "d4eb6e29-4a72-4046-9332-9f9e7f4f6245.4c6490aceeada19c7b905ed6c0f1534c77a08004.$langSuffix.xlf"
    $ws.Range("D4").Value2 = $targetFile4
    $ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.$langSuffix/ci/ht/$targetFile4", "", "", $targetFile4) | Out-Null
    $ws.Range("E4").Value2 = $handoffDatetime
    $ws.Range("E4").NumberFormat = $dateFmt
    $ws.Range("H4").Value2 = "0001-01-01 00:00:00"
    $ws.Range("H4").NumberFormat = $dateFmt
    $ws.Range("J4").Value2 = "Include"

    # Row 5 - f1d821ef-a8ce-4150-9380-8d83689f79af
    $ws.Range("A5").Value2 = "f1d821ef-a8ce-4150-9380-8d83689f79af.md"
    $ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest.$langSuffix/blob/master/e2e/f1d821ef-a8ce-4150-9380-8d83689f79af.md", "", "", "f1d821ef-a8ce-4150-9380-8d83689f79af.md") | Out-Null
    $ws.Range("B5").Value2 = ".md"
    $ws.Range("C5").Value2 = "Ready for handoff"
    $targetFile5 = "f1d821ef-a8ce-4150-9380-8d83689f79af.b387254a5552561c31c6576c1aca265deac59e71.$langSuffix.xlf"
    $ws.Range("D5").Value2 = $targetFile5
    $ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.$langSuffix/ci/ht/$targetFile5", "", "", $targetFile5) | Out-Null
    $ws.Range("E5").Value2 = $handoffDatetime
    $ws.Range("E5").NumberFormat = $dateFmt
    $ws.Range("H5").Value2 = "0001-01-01 00:00:00"
    $ws.Range("H5").NumberFormat = $dateFmt
    $ws.Range("J5").Value2 = "Include"
}

Add-DetailRows $zhcn "zh-cn" "2016-03-21 14:37:54"
Add-DetailRows $dede "de-de" "2016-03-21 14:37:58"
